# Append a new data row (2025/12/31, 逃离鸭科夫, 1121) to the mod-count sheet,
# matching the style/formatting of the preceding rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Find the last used row in column A (the Date column) and the row right
# after it, which is where the new record goes.
$xlUp = [Microsoft.Office.Interop.Excel.XlDirection]::xlUp
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End($xlUp).Row
$newRow = $lastRow + 1

# Format the date cell as Text first so Excel stores the literal
# "2025/12/31" string instead of auto-converting it to a date serial value
# (matches how the existing rows store their Date column as text).
$ws.Cells.Item($newRow, 1).NumberFormat = "@"

$ws.Cells.Item($newRow, 1).Value = "2025/12/31"
$ws.Cells.Item($newRow, 2).Value = "逃离鸭科夫"
$ws.Cells.Item($newRow, 3).Value = 1121

# Copy the formatting (centered alignment / style) from the row above so the
# new row visually matches the rest of the table.
$srcRange = $ws.Range($ws.Cells.Item($lastRow, 1), $ws.Cells.Item($lastRow, 3))
$dstRange = $ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, 3))
$srcRange.Copy()
$dstRange.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false
